# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - The second worker "mora" record (period 2506) is removed entirely
#   (the whole row 17 is deleted, shifting the closing signature rows up).
# - The remaining record's period changes from 2507 -> 2508.
# - VALOR MORA total (E11) now reflects only the single remaining record.
# - Cant. Periodos (F13) drops from 2 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for the second ("2506") worker record - the sheet keeps
# only one detail row going forward.
$ws.Rows("17").Delete()

# Update the period code of the remaining detail row.
$ws.Range("E16").Value = "2508"

# Update the dependent summary figures.
$ws.Range("E11").Value = 56940
$ws.Range("F13").Value = 1
